# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(1,5,13,8,9,4,5,3,2,3,3,5,2,2,2,0,2,1,2,4,1,2,2,2,2,1,0,1,0,2,3,3,4,3,3,0,0,0,3,1,3,2,1,4,2,1,3,1,3,2,2,1,1,1,4)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = $newK[$i]
}
